$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65 (pushes old row 65.. down to 66..), copying
# formatting from the row above (row 64), matching the "law-why" tag
# block getting a new FAQ pair.
$ws.Rows.Item(65).Insert()

# Row 64 previously duplicated "ทำไมต้องเก็บภาษี" (same text as row 60);
# replace it with the new question text.
$ws.Range("B64").Value = "เก็บภาษีทำไม"

# New row 65: continuation of the "law-why" tag group.
$ws.Range("A65").Value = "law-why"
$ws.Range("B65").Value = "เก็บภาษีเพื่อ"

# Reflect the author's final selection/viewport on the sheet.
$ws.Range("G68").Select()
